$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to store a literal text value even when the
# string looks like a number (e.g. "1.001", "21.20"), without leaving
# a residual custom number format behind on the cell.
function Set-TextValue($rangeRef, $text) {
    $rng = $ws.Range($rangeRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '27.160.16'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '1.904.87'
$ws.Range('E3').Value = '  +0.77%  '
Set-TextValue 'D4' '1.001'
$ws.Range('E4').Value = '  +0.03%  '
Set-TextValue 'D5' '306.28'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('E6').Value = '  -0.01%  '
Set-TextValue 'D7' '0.5235'
$ws.Range('E7').Value = '  +1.70%  '
Set-TextValue 'D8' '0.3767'
$ws.Range('E8').Value = '  -0.01%  '
Set-TextValue 'D9' '0.07247'
$ws.Range('E9').Value = '  +0.47%  '
Set-TextValue 'D10' '21.20'
$ws.Range('E10').Value = '  +0.12%  '
Set-TextValue 'D11' '0.9025'
$ws.Range('E11').Value = '  -0.18%  '
Set-TextValue 'D12' '0.08524'
$ws.Range('E12').Value = '  +11.41%  '
$ws.Range('D13').Value = '1.917.29'
$ws.Range('E13').Value = '  +1.40%  '
Set-TextValue 'D14' '96.97'
$ws.Range('E14').Value = '  +2.25%  '
Set-TextValue 'D15' '5.294'
$ws.Range('E15').Value = '  +0.39%  '
Set-TextValue 'D16' '1.001'
$ws.Range('E16').Value = '  -0.27%  '
Set-TextValue 'D17' '0.000008634'
$ws.Range('E17').Value = '  +1.88%  '
$ws.Range('E18').Value = '  +0.87%  '
Set-TextValue 'D19' '1.000'
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').Value = '27.194.99'
$ws.Range('E20').Value = '  +0.34%  '
Set-TextValue 'D21' '5.071'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').Value = '2.158.29'
$ws.Range('E22').Value = '  +1.49%  '
$ws.Range('E23').Value = '  +0.50%  '
Set-TextValue 'D24' '6.439'
$ws.Range('E24').Value = '  +0.58%  '
$ws.Range('E25').Value = '  +2.37%  '
Set-TextValue 'D26' '147.17'
$ws.Range('E26').Value = '  +0.95%  '
Set-TextValue 'D27' '18.25'
$ws.Range('E27').Value = '  +0.91%  '
Set-TextValue 'D28' '1.747'
$ws.Range('E28').Value = '  -1.91%  '
Set-TextValue 'D29' '114.96'
$ws.Range('E29').Value = '  +0.43%  '
Set-TextValue 'D30' '4.924'
$ws.Range('E30').Value = '  -0.50%  '
Set-TextValue 'D31' '4.819'
$ws.Range('E31').Value = '  -0.19%  '
Set-TextValue 'D32' '0.09283'
$ws.Range('E32').Value = '  +1.14%  '
Set-TextValue 'D33' '0.8066'
$ws.Range('E33').Value = '  +3.04%  '
Set-TextValue 'D34' '0.05054'
$ws.Range('E34').Value = '  -0.59%  '
Set-TextValue 'D35' '1.244'
$ws.Range('E35').Value = '  +0.67%  '
Set-TextValue 'D36' '3.453'
$ws.Range('E36').Value = '  +5.03%  '
Set-TextValue 'D37' '2.952'
$ws.Range('E37').Value = '  -1.68%  '
Set-TextValue 'D38' '2.618'
$ws.Range('E38').Value = '  -0.43%  '
Set-TextValue 'D39' '0.5715'
$ws.Range('E39').Value = '  +2.37%  '
Set-TextValue 'D40' '0.02001'
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('E41').Value = '  -0.10%  '
Set-TextValue 'D42' '9.148'
$ws.Range('E42').Value = '  +1.69%  '
Set-TextValue 'D43' '6.641'
$ws.Range('E43').Value = '  +0.23%  '
Set-TextValue 'D44' '116.10'
$ws.Range('E44').Value = '  -1.33%  '
$ws.Range('E45').Value = '  +0.45%  '
Set-TextValue 'D46' '0.4870'
$ws.Range('E46').Value = '  +1.48%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D47' '10.21'
$ws.Range('E47').Value = '  -0.18%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D48' '0.9997'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('E49').Value = '  +1.38%  '
Set-TextValue 'D50' '37.55'
$ws.Range('E50').Value = '  -0.21%  '
Set-TextValue 'D51' '64.27'
$ws.Range('E51').Value = '  +0.37%  '
